$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.934.12'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.033.22'
$ws.Range('E3').Value = '  -0.92%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.22'
$ws.Range('E5').Value = '  -0.62%  '

$ws.Range('E6').Value = '  -0.63%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.71'
$ws.Range('E7').Value = '  +3.31%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.382'
$ws.Range('E9').Value = '  -0.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0817'
$ws.Range('E10').Value = '  +1.07%  '

$ws.Range('E11').Value = '  +0.25%  '

$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.57'
$ws.Range('E12').Value = '  -0.55%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.334.32'
$ws.Range('E13').Value = '  -0.90%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.37'
$ws.Range('E14').Value = '  +2.39%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.762'
$ws.Range('E15').Value = '  +1.42%  '

$ws.Range('E16').Value = '  -2.35%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.045.66'
$ws.Range('E17').Value = '  -0.43%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.842.76'
$ws.Range('E18').Value = '  -0.31%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.86'
$ws.Range('E19').Value = '  +0.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.95'
$ws.Range('E20').Value = '  -5.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.48'
$ws.Range('E22').Value = '  -0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  -1.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  +0.57%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.11'
$ws.Range('E26').Value = '  +0.50%  '

$ws.Range('E27').Value = '  +0.12%  '

$ws.Range('E28').Value = '  -4.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.92'
$ws.Range('E29').Value = '  -0.49%  '

$ws.Range('E30').Value = '  -3.73%  '

$ws.Range('E31').Value = '  +0.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.16'
$ws.Range('E32').Value = '  +4.84%  '

$ws.Range('E33').Value = '  -2.41%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0607'
$ws.Range('E34').Value = '  -0.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.51'
$ws.Range('E35').Value = '  -1.80%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.46'
$ws.Range('E36').Value = '  +5.96%  '

$ws.Range('E37').Value = '  -2.00%  '

$ws.Range('E38').Value = '  -0.72%  '

$ws.Range('E39').Value = '  +0.10%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.525.45'
$ws.Range('E40').Value = '  +2.60%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.24'
$ws.Range('E41').Value = '  +4.21%  '

$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.19'
$ws.Range('E43').Value = '  -0.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.84'
$ws.Range('E44').Value = '  -1.55%  '

$ws.Range('E45').Value = '  -1.16%  '

$ws.Range('E46').Value = '  -1.72%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.95'
$ws.Range('E47').Value = '  -4.25%  '

$ws.Range('E48').Value = '  -0.60%  '

$ws.Range('E49').Value = '  -0.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.10'
$ws.Range('E50').Value = '  +0.44%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.223.10'
$ws.Range('E51').Value = '  -0.92%  '
